$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OCT-2020")

# Insert 3 fresh blank rows right before the legend block (old rows 32:34),
# pushing the legend (old rows 32-36) down to new rows 35-39 while rows
# 25-31 (already blank/unused inside the old A1:G36 dimension) become the
# new data rows.
$ws.Rows("32:34").Insert()

# ---- Row 25 (24 Oct 2020 - Week off) : format like row 18 ----
$ws.Range("A18:G18").Copy()
$ws.Range("A25:G25").PasteSpecial(-4122)
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 44128
$ws.Cells.Item(25, 4).Value = "Week off"

# ---- Row 26 (25 Oct 2020 - Week off) : format like row 18 ----
$ws.Range("A18:G18").Copy()
$ws.Range("A26:G26").PasteSpecial(-4122)
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = 44129
$ws.Cells.Item(26, 4).Value = "Week off"

# ---- Row 27 (26 Oct 2020 - Holiday) : format like row 19 ----
$ws.Range("A19:G19").Copy()
$ws.Range("A27:G27").PasteSpecial(-4122)
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = 44130
$ws.Cells.Item(27, 4).Value = "Holiday"

# ---- Row 28 (27 Oct 2020 - QMVAR work) : format like row 22, col A like E18 ----
$ws.Range("A22:G22").Copy()
$ws.Range("A28:G28").PasteSpecial(-4122)
$ws.Range("E18").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 44131
$ws.Cells.Item(28, 3).Value = "QMVAR"
$ws.Cells.Item(28, 4).Value = "QMVAR - alignment issues"
$ws.Cells.Item(28, 5).Value = 0.7
$ws.Cells.Item(28, 6).Value = "Completed"

# ---- Row 29 (28 Oct 2020 - QMVAR work) : format like row 22, col A like E18 ----
$ws.Range("A22:G22").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)
$ws.Range("E18").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 44132
$ws.Cells.Item(29, 3).Value = "QMVAR"
$ws.Cells.Item(29, 4).Value = "QMVAR -Alignment issues fixing"
$ws.Cells.Item(29, 5).Value = 0.8
$ws.Cells.Item(29, 6).Value = "Completed"

# ---- Row 30 (29 Oct 2020 - placeholder, no task yet) : format like row 22, col A like E18 ----
$ws.Range("A22:G22").Copy()
$ws.Range("A30:G30").PasteSpecial(-4122)
$ws.Range("E18").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 44133

# ---- Row 31 (30 Oct 2020 - placeholder, no task yet) : format like row 22, col A like E18 ----
$ws.Range("A22:G22").Copy()
$ws.Range("A31:G31").PasteSpecial(-4122)
$ws.Range("E18").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 44134

$excel.CutCopyMode = $false

# ---- Update the view: scroll position + active selection ----
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D33").Select()
